# Update "想去人数" (want-to-go count) values in the F column for the
# 展览 (Exhibitions) and 全部类型 (All types) sheets, per the commit's
# regenerated data output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value, keyed by sheet because "全部类型" has one extra row
# (a cancelled/zero-interest event) inserted before the matching rows.
$exhibitUpdates = @{
    5  = 5
    6  = 543
    7  = 7613
    8  = 481
    9  = 197
    10 = 1075
    11 = 621
    12 = 7
    13 = 26
    14 = 172
    15 = 3
    16 = 199
    17 = 743
}

$allTypesUpdates = @{
    5  = 5
    7  = 543
    8  = 7613
    9  = 481
    10 = 197
    11 = 1075
    12 = 621
    13 = 7
    14 = 26
    15 = 172
    16 = 3
    17 = 199
    18 = 743
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}
